$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A
$ws.Range("A1").Value = "Impact"
$ws.Range("A2").Value = "Claps"
$ws.Range("A3").Value = "Nono"
$ws.Range("A4").Value = "Goobie"
$ws.Range("A5").Value = "Stastro"
$ws.Range("A6").Value = "Opps"
$ws.Range("A7").Value = "Zombie"

# Column B
$ws.Range("B1").Value = "Crew"
$ws.Range("B2").Value = "Padge"
$ws.Range("B3").Value = "Hoops"
$ws.Range("B4").Value = "Pizza"
$ws.Range("B5").Value = "Ups"
$ws.Range("B6").Value = "GoGo"
$ws.Range("B7").Value = "Nom Nom"

# Column C
$ws.Range("C1").Value = "Cove"
$ws.Range("C2").Value = "Hawma"
$ws.Range("C3").Value = "G-Poppy"
$ws.Range("C4").Value = "Unicorn"

# Column D
$ws.Range("D1").Value = "Workcrew"
$ws.Range("D2").Value = "Laddy"
$ws.Range("D3").Value = "Burning Bush"
$ws.Range("D4").Value = "Bolt"
$ws.Range("D5").Value = "Bonez"
$ws.Range("D6").Value = "Sea"

# Column E (new)
$ws.Range("E1").Value = "P-Staff"
$ws.Range("E2").Value = "Indi"
$ws.Range("E3").Value = "Screams"
$ws.Range("E4").Value = "Bow Wow"
$ws.Range("E5").Value = "Tross"

# Header row formatting extends into F1 (bold, no value)
$ws.Range("E1:F1").Font.Bold = $true

$ws.Range("E6").Select()
